$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 28: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I28").Value = "sv"
$ws.Range("J28").Value = "Statement-opinion"

# Row 29: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I29").Value = "sv"
$ws.Range("J29").Value = "Statement-opinion"

# Row 32: sv/Statement-opinion -> sd/Statement-non-opinion
$ws.Range("I32").Value = "sd"
$ws.Range("J32").Value = "Statement-non-opinion"

# Row 45: qy/Yes-No-Question -> ba/Appreciation
$ws.Range("I45").Value = "ba"
$ws.Range("J45").Value = "Appreciation"

# Row 50: sd/Statement-non-opinion -> qy/Yes-No-Question
$ws.Range("I50").Value = "qy"
$ws.Range("J50").Value = "Yes-No-Question"
